$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NSL-KDD-Normalisation")
$ws.Range("J2").Value2 = 0.9500510136184181
$ws.Range("K2").Value2 = 0.96
$ws.Range("L2").Value2 = 0.86100326137599
$ws.Range("M2").Value2 = 0.9078107090224332
$ws.Range("J3").Value2 = 0.9406467639622056
$ws.Range("K3").Value2 = 0.9586405322783672
$ws.Range("L3").Value2 = 0.8279235906196615
$ws.Range("M3").Value2 = 0.8885
$ws.Range("J4").Value2 = 0.9417557556669476
$ws.Range("K4").Value2 = 0.9591544249372984
$ws.Range("L4").Value2 = 0.8314955738468707
$ws.Range("M4").Value2 = 0.890774477996839
$ws.Range("J5").Value2 = 0.9435301423945349
$ws.Range("K5").Value2 = 0.948905109489051
$ws.Range("L5").Value2 = 0.8479577574157477
$ws.Range("M5").Value2 = 0.8955958336750595
$ws.Range("J6").Value2 = 0.81031805882092
$ws.Range("K6").Value2 = 0.6179775280898876
$ws.Range("M6").Value2 = 0.7260028194284249

$ws = $wb.Worksheets.Item("NSL-KDD-RF feature s")
$ws.Range("J2").Value2 = 0.9440624584128111
$ws.Range("K2").Value2 = 0.9383677615983745
$ws.Range("L2").Value2 = 0.8606926541388414
$ws.Range("M2").Value2 = 0.8978533819360064
$ws.Range("J3").Value2 = 0.9424655103579825
$ws.Range("K3").Value2 = 0.958943234559086
$ws.Range("L3").Value2 = 0.8342910389812083
$ws.Range("M3").Value2 = 0.8922846939622954
$ws.Range("J4").Value2 = 0.9430421860444483
$ws.Range("K4").Value2 = 0.9592018528416176
$ws.Range("L4").Value2 = 0.8361546824041001
$ws.Range("M4").Value2 = 0.8934616661135081
$ws.Range("J5").Value2 = 0.9114137426252051
$ws.Range("K5").Value2 = 0.9000360230547551
$ws.Range("L5").Value2 = 0.7760521820158409
$ws.Range("M5").Value2 = 0.8334584271536986
$ws.Range("J6").Value2 = 0.8760590870780287
$ws.Range("K6").Value2 = 0.8130905342724618
$ws.Range("L6").Value2 = 0.7350520267122224
$ws.Range("M6").Value2 = 0.7721044045676998

$ws = $wb.Worksheets.Item("NSL-KDD-RFE feature s")
$ws.Range("J2").Value2 = 0.9503171716275562
$ws.Range("K2").Value2 = 0.9608386761393173
$ws.Range("L2").Value2 = 0.8611585649945643
$ws.Range("M2").Value2 = 0.9082719082719082
$ws.Range("J3").Value2 = 0.9406024042940159
$ws.Range("K3").Value2 = 0.9584681769147788
$ws.Range("L3").Value2 = 0.8279235906196615
$ws.Range("M3").Value2 = 0.8884259645029581
$ws.Range("J4").Value2 = 0.9405580446258262
$ws.Range("K4").Value2 = 0.9587907144142523
$ws.Range("L4").Value2 = 0.8274576797639385
$ws.Range("M4").Value2 = 0.8882960986995665
$ws.Range("J5").Value2 = 0.9454819677948809
$ws.Range("K5").Value2 = 0.9579817158931083
$ws.Range("L5").Value2 = 0.8462494176114304
$ws.Range("M5").Value2 = 0.8986558918116599
$ws.Range("J6").Value2 = 0.9023200106463204
$ws.Range("K6").Value2 = 0.8710807496934665
$ws.Range("L6").Value2 = 0.7723248951700574
$ws.Range("M6").Value2 = 0.8187355943365163

$ws = $wb.Worksheets.Item("NSL-KDD-Linear reg feature s")
$ws.Range("J2").Value2 = 0.9471676351860888
$ws.Range("K2").Value2 = 0.9555555555555556
$ws.Range("L2").Value2 = 0.8547911166330175
$ws.Range("M2").Value2 = 0.9023690466431674
$ws.Range("J3").Value2 = 0.9330169010335803
$ws.Range("K3").Value2 = 0.9597090095131505
$ws.Range("L3").Value2 = 0.7990371175648393
$ws.Range("M3").Value2 = 0.8720338983050847
$ws.Range("J4").Value2 = 0.9342589717428914
$ws.Range("K4").Value2 = 0.9599183521989237
$ws.Range("L4").Value2 = 0.80338561888492
$ws.Range("M4").Value2 = 0.87470409198512
$ws.Range("J5").Value2 = 0.939759570598412
$ws.Range("K5").Value2 = 0.939152981849611
$ws.Range("L5").Value2 = 0.8437645597142414
$ws.Range("M5").Value2 = 0.8889070680628273
$ws.Range("J6").Value2 = 0.8780109124783747
$ws.Range("K6").Value2 = 0.7769103738177451
$ws.Range("L6").Value2 = 0.8036962261220686
$ws.Range("M6").Value2 = 0.7900763358778624
